$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Version & History" - add a new V2.4 history row (row 19)
# ---------------------------------------------------------------------------
$wsHist = $wb.Worksheets.Item(1)

# Copy formatting (styles) of the previous history row (18) into the new row (19)
$wsHist.Range("A18:E18").Copy($wsHist.Range("A19:E19"))
$wsHist.Rows("19").RowHeight = 83.95

# Column A: Version
$wsHist.Range("A19").Value = "V2.4"

# Column B: rich-text description of the change ("Changed: ..." + details in a
# different font size, matching the formatting pattern used by the previous
# "Added:" row).
$changeText = "Changed:`nLongitudinal RVX → Radar Longitudinal Relative Velocity`nLongitudinal EGO → Radar Longitudinal Distance From EGO`nLateral RVY → Radar Lateral Relative Velocity`nLateral EGO → Radar Lateral Distance From EGO`n"
$wsHist.Range("B19").Value = $changeText

$headLen = "Changed:`n".Length
$totalLen = $changeText.Length

$headRange = $wsHist.Range("B19").Characters(1, $headLen)
$headRange.Font.Size = 11
$headRange.Font.Color = 0
$headRange.Font.Name = "Calibri"

$bodyRange = $wsHist.Range("B19").Characters($headLen + 1, $totalLen - $headLen)
$bodyRange.Font.Size = 12
$bodyRange.Font.Color = 0
$bodyRange.Font.Name = "Calibri"

# Column C: Author
$wsHist.Range("C19").Value = "Zborai Attila"

# Columns D (Date) and E (Status) keep the same values as row 18
# (04.11.2017 / Draft version), already copied above.

# ---------------------------------------------------------------------------
# Sheet 2: "CommunicationMatrix" - rename the Radar Sensor signals
# ---------------------------------------------------------------------------
$wsMatrix = $wb.Worksheets.Item(2)

$wsMatrix.Range("B15").Value = "Radar Longitudinal Relative Velocity"
$wsMatrix.Range("B16").Value = "Radar Longitudinal Distance From EGO"
$wsMatrix.Range("B17").Value = "Radar Lateral Relative Velocity"
$wsMatrix.Range("B18").Value = "Radar Lateral Distance From EGO"

# Widen column B to fit the longer signal names
$wsMatrix.Columns("B").ColumnWidth = 23.9

# Update the selected/active cell on each sheet to match the saved view state
$wsHist.Range("B26").Select()

$wsMatrix.Activate()
$wsMatrix.Range("B15").Select()
